$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("STREAMS")

# Insert a new row at position 15 (shifts existing rows 15-22 down to 16-23)
$ws.Rows.Item(15).Insert()

$ws.Range("A15").Value = "Liquid molar fraction water"
$ws.Range("B15").Value = "xH2O"
$cols = @("C","D","E","F","G","H","I","J","K","L","M")
foreach ($col in $cols) {
    $ws.Range($col + "15").Value = 0
}

# Append a new row 24 at the end (after old row 22, now shifted to row 23)
$ws.Range("A24").Value = "Gas molar fraction water"
$ws.Range("A24").Font.Bold = $true
$ws.Range("A24").Font.Color = 0
$ws.Range("B24").Value = "yH2O"
foreach ($col in $cols) {
    $ws.Range($col + "24").Value = 0
}

# Update pressure row values (K6, L6 from 500000 to 100000) with scientific format
$ws.Range("K6").Value = 100000
$ws.Range("L6").Value = 100000
$ws.Range("K6:L6").NumberFormat = "0.00E+00"

# Selections on each sheet
$wsCompounds = $wb.Worksheets.Item("COMPOUNDS")
[void]$wsCompounds.Range("H11").Select()

$wsUnitOps = $wb.Worksheets.Item("UNIT OPERATIONS")
[void]$wsUnitOps.Range("I1").Select()

[void]$ws.Range("B25").Select()
[void]$ws.Activate()
